# Apply the new default table style to the three data tables (slides 14-16)
# that previously used the bare "Table_0" style defined in tableStyles.xml.
$p = $ppt.ActivePresentation

$newStyleId = "{4C35CE88-87D2-4678-85A5-9413F435D433}"
$tableSlideIndexes = 14, 15, 16

foreach ($idx in $tableSlideIndexes) {
    $slide = $p.Slides.Item($idx)
    foreach ($shp in $slide.Shapes) {
        if ($shp.HasTable) {
            $shp.Table.ApplyStyle($newStyleId)
        }
    }
}
